{"js": "// Apply Papiamento translation fixes to pap_A-E-C-P.docx\n// Each entry: unique search text (old) -> replacement text (new)\nconst replacements = [\n  {\n    find: \"skucha aktivo, fasilitashon den kolaborashon, krea konsientisashon\",\n    replace: \"skucha aktivo, enbolb\u00ed mayornan aktivamente, krea konsientisashon\"\n  },\n  {\n    find: \"Eh\u00e8mpel di un Pr\u00e1ktika di Grupo (praktik\u00e1 Tempu pa abo ku bo yu):  \",\n    replace: \"Eh\u00e8mpel di un Pr\u00e1ktika di Grupo (praktik\u00e1 Tempu-pa Abo-ku Bo yu):  \"\n  },\n  {\n    find: \"Abo: Mi por komprond\u00e9 ku bo lo sinti bo ink\u00f3modo den kuminsamentu\",\n    replace: \"Abo: Mi por komprond\u00e9 ku lo bo por sinti ink\u00f3modo den kuminsamentu\"\n  },\n  {\n    find: \"E ta sosten\u00e9 kon solushon\u00e1 problema i ta yuda mayornan refleh\u00e1\",\n    replace: \"E ta sosten\u00e9 e abilidat di solushon\u00e1 problema i ta yuda mayornan refleh\u00e1\"\n  },\n  {\n    find: \"riba pr\u00e1ktika di Tempu Abo ku bo Yu:  \",\n    replace: \"riba pr\u00e1ktika di Tempu-pa Abo-ku Bo yu:  \"\n  },\n  {\n    find: \"Resumen: Dor di Pasa Tempu Abo ku Bo Yu i laga nan tuma e liderazgo\",\n    replace: \"Resumen: Dor di Pasa Tempu-pa Abo-ku Bo yu i laga nan tuma e liderazgo\"\n  },\n  {\n    find: \"e seshon tokante Tempu pa abo ku bo Yu.) \",\n    replace: \"e seshon tokante Tempu-pa Abo-ku Bo yu.) \"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + find);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply Papiamento translation fixes to pap_A-E-C-P.docx\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"skucha aktivo, fasilitashon den kolaborashon, krea konsientisashon\"; Replace = \"skucha aktivo, enbolb\u00ed mayornan aktivamente, krea konsientisashon\" },\n    @{ Find = \"Eh\u00e8mpel di un Pr\u00e1ktika di Grupo (praktik\u00e1 Tempu pa abo ku bo yu):  \"; Replace = \"Eh\u00e8mpel di un Pr\u00e1ktika di Grupo (praktik\u00e1 Tempu-pa Abo-ku Bo yu):  \" },\n    @{ Find = \"Abo: Mi por komprond\u00e9 ku bo lo sinti bo ink\u00f3modo den kuminsamentu\"; Replace = \"Abo: Mi por komprond\u00e9 ku lo bo por sinti ink\u00f3modo den kuminsamentu\" },\n    @{ Find = \"E ta sosten\u00e9 kon solushon\u00e1 problema i ta yuda mayornan refleh\u00e1\"; Replace = \"E ta sosten\u00e9 e abilidat di solushon\u00e1 problema i ta yuda mayornan refleh\u00e1\" },\n    @{ Find = \"riba pr\u00e1ktika di Tempu Abo ku bo Yu:  \"; Replace = \"riba pr\u00e1ktika di Tempu-pa Abo-ku Bo yu:  \" },\n    @{ Find = \"Resumen: Dor di Pasa Tempu Abo ku Bo Yu i laga nan tuma e liderazgo\"; Replace = \"Resumen: Dor di Pasa Tempu-pa Abo-ku Bo yu i laga nan tuma e liderazgo\" },\n    @{ Find = \"e seshon tokante Tempu pa abo ku bo Yu.) \"; Replace = \"e seshon tokante Tempu-pa Abo-ku Bo yu.) \" }\n)\n\nforeach ($item in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $found = $find.Execute($item.Find, $false, $false, $false, $false, $false, $true, 1, $false, $item.Replace, 2)\n    if (-not $found) {\n        throw \"Could not find text: $($item.Find)\"\n    }\n}\n"}
